$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.158.18"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.36"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.27"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4968"
$ws.Range("E7").Value = "  -3.44%  "

$ws.Range("E8").Value = "  -1.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09821"
$ws.Range("E9").Value = "  +24.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.26"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.447"
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.817.61"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.318"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001140"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.68"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06667"
$ws.Range("E19").Value = "  +1.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.024"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.210.56"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.246"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.93"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.78"
$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.032.65"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.427"
$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.13"
$ws.Range("E30").Value = "  -1.58%  "

$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.615"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06746"
$ws.Range("E35").Value = "  -6.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.020"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.978"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6233"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("E42").Value = "  +2.14%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5946"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.275"
$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.22"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.952"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06794"
$ws.Range("E51").Value = "  -0.95%  "
